$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value()
    $dVal = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
}
